# Inserts a new weekly price record as row 66 (Fecha 2021-12-23, serial 44553),
# pushing the existing rows 66-101 down to 67-102.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("66:66").Insert()

$ws.Range("A66").Value = 2
$ws.Range("B66").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C66").Value = "Coquimbo"
$ws.Range("D66").Value = 44553
$ws.Range("E66").Value = 4
$ws.Range("F66").Value = 100112024
$ws.Range("G66").Value = "Choclo"
$ws.Range("H66").Value = "Dulce o Americano"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 600
$ws.Range("K66").Value = 13000
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = 14000
$ws.Range("N66").Value = "$/malla 70 unidades"
$ws.Range("O66").Value = "Provincia de Limarí"
$ws.Range("P66").Value = 200
$ws.Range("Q66").Value = 70
$ws.Range("R66").Value = "Hortaliza"
